$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (40 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1142.3636
$ws.Range("I28").Value = 1393
$ws.Range("J28").Value = 641.0909
$ws.Range("K28").Value = 1393
$ws.Range("L28").Value = 641.0909
$ws.Range("M28").Value = -908
$ws.Range("N28").Value = -1611.0909
$ws.Range("H38").Value = 73.416664
$ws.Range("I38").Value = 43.727272
$ws.Range("J38").Value = 400
$ws.Range("K38").Value = 131.181816
$ws.Range("L38").Value = 1200
$ws.Range("M38").Value = 240.818184
$ws.Range("N38").Value = -1944
$ws.Range("H40").Value = 1176.6
$ws.Range("I40").Value = 1127.2222
$ws.Range("J40").Value = 1250.6666
$ws.Range("K40").Value = 1127.2222
$ws.Range("L40").Value = 1250.6666
$ws.Range("M40").Value = -952.2221999999999
$ws.Range("N40").Value = -1600.6666
$ws.Range("H96").Value = 907.6667
$ws.Range("I96").Value = 883.6
$ws.Range("K96").Value = 2650.8
$ws.Range("M96").Value = -1277.8
$ws.Range("H106").Value = 1382
$ws.Range("I106").Value = 962.3333
$ws.Range("K106").Value = 962.3333
$ws.Range("M106").Value = -331.3333
$ws.Range("H133").Value = 68475
$ws.Range("J133").Value = 68475
$ws.Range("L133").Value = 68475
$ws.Range("N133").Value = -78595
$ws.Range("H137").Value = 1360
$ws.Range("I137").Value = 1360
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 4080
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -1530
$ws.Range("N137").ClearContents()

# --- Sheet: ARM (37 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 505.57144
$ws.Range("I5").Value = 1024.75
$ws.Range("J5").Value = 116.1875
$ws.Range("K5").Value = 1024.75
$ws.Range("L5").Value = 116.1875
$ws.Range("M5").Value = -912.75
$ws.Range("N5").Value = -340.1875
$ws.Range("H45").Value = 14891.35
$ws.Range("I45").Value = 12306.6
$ws.Range("J45").Value = 22645.6
$ws.Range("K45").Value = 12306.6
$ws.Range("L45").Value = 22645.6
$ws.Range("M45").Value = -11929.6
$ws.Range("N45").Value = -23399.6
$ws.Range("H74").Value = 4966.75
$ws.Range("I74").Value = 585.25
$ws.Range("K74").Value = 585.25
$ws.Range("M74").Value = 288.75
$ws.Range("H77").Value = 4966.75
$ws.Range("I77").Value = 585.25
$ws.Range("K77").Value = 2926.25
$ws.Range("M77").Value = 1441.75
$ws.Range("H110").Value = 20976.8
$ws.Range("I110").Value = 33269
$ws.Range("J110").Value = 6928.5713
$ws.Range("K110").Value = 33269
$ws.Range("L110").Value = 6928.5713
$ws.Range("M110").Value = -31224
$ws.Range("N110").Value = -11018.5713
$ws.Range("H122").Value = 2536
$ws.Range("I122").Value = 2372
$ws.Range("K122").Value = 7116
$ws.Range("M122").Value = -4666
$ws.Range("H132").Value = 1865.3549
$ws.Range("I132").Value = 1649.1724
$ws.Range("K132").Value = 4947.5172
$ws.Range("M132").Value = -2417.5172

# --- Sheet: BSM (34 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 505.57144
$ws.Range("I4").Value = 1024.75
$ws.Range("J4").Value = 116.1875
$ws.Range("K4").Value = 1024.75
$ws.Range("L4").Value = 116.1875
$ws.Range("M4").Value = -909.75
$ws.Range("N4").Value = -346.1875
$ws.Range("H36").Value = 21447.6
$ws.Range("I36").Value = 1809.5
$ws.Range("K36").Value = 1809.5
$ws.Range("M36").Value = -1275.5
$ws.Range("H64").Value = 1539.6
$ws.Range("J64").Value = 1461
$ws.Range("L64").Value = 1461
$ws.Range("N64").Value = -1911
$ws.Range("H67").Value = 1539.6
$ws.Range("J67").Value = 1461
$ws.Range("L67").Value = 1461
$ws.Range("N67").Value = -3021
$ws.Range("H86").Value = 1916.5
$ws.Range("I86").Value = 1916.5
$ws.Range("K86").Value = 1916.5
$ws.Range("M86").Value = -793.5
$ws.Range("H89").Value = 1916.5
$ws.Range("I89").Value = 1916.5
$ws.Range("K89").Value = 9582.5
$ws.Range("M89").Value = -3966.5
$ws.Range("H134").Value = 2832.1904
$ws.Range("I134").Value = 2505.9333
$ws.Range("J134").Value = 3647.8333
$ws.Range("K134").Value = 7517.7999
$ws.Range("L134").Value = 10943.4999
$ws.Range("M134").Value = -4982.7999
$ws.Range("N134").Value = -16013.4999

# --- Sheet: CRP (14 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 54915.95
$ws.Range("I31").Value = 64782.75
$ws.Range("J31").Value = 15448.75
$ws.Range("K31").Value = 64782.75
$ws.Range("L31").Value = 15448.75
$ws.Range("M31").Value = -64487.75
$ws.Range("N31").Value = -16038.75
$ws.Range("H34").Value = 54915.95
$ws.Range("I34").Value = 64782.75
$ws.Range("J34").Value = 15448.75
$ws.Range("K34").Value = 64782.75
$ws.Range("L34").Value = 15448.75
$ws.Range("M34").Value = -64580.75
$ws.Range("N34").Value = -15852.75

# --- Sheet: CUL (30 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 114
$ws.Range("J12").Value = 131.66667
$ws.Range("L12").Value = 395.00001
$ws.Range("N12").Value = -741.00001
$ws.Range("H86").Value = 618.1429000000001
$ws.Range("I86").Value = 361
$ws.Range("J86").Value = 961
$ws.Range("K86").Value = 1083
$ws.Range("L86").Value = 2883
$ws.Range("M86").Value = 103
$ws.Range("N86").Value = -5255
$ws.Range("H87").Value = 8318.666999999999
$ws.Range("I87").Value = 5608.5
$ws.Range("K87").Value = 16825.5
$ws.Range("M87").Value = -15577.5
$ws.Range("H89").Value = 618.1429000000001
$ws.Range("I89").Value = 361
$ws.Range("J89").Value = 961
$ws.Range("K89").Value = 3249
$ws.Range("L89").Value = 8649
$ws.Range("M89").Value = 2679
$ws.Range("N89").Value = -20505
$ws.Range("H90").Value = 8318.666999999999
$ws.Range("I90").Value = 5608.5
$ws.Range("K90").Value = 50476.5
$ws.Range("M90").Value = -44236.5
$ws.Range("H133").Value = 20624.875
$ws.Range("I133").Value = 23999.8
$ws.Range("K133").Value = 71999.39999999999
$ws.Range("M133").Value = -66939.39999999999

# --- Sheet: LTW (65 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 15374.75
$ws.Range("I7").Value = 18499.666
$ws.Range("J7").Value = 6000
$ws.Range("K7").Value = 18499.666
$ws.Range("L7").Value = 6000
$ws.Range("M7").Value = -18387.666
$ws.Range("N7").Value = -6224
$ws.Range("H22").Value = 1217.4286
$ws.Range("I22").Value = 982.625
$ws.Range("J22").Value = 1530.5
$ws.Range("K22").Value = 982.625
$ws.Range("L22").Value = 1530.5
$ws.Range("M22").Value = -687.625
$ws.Range("N22").Value = -2120.5
$ws.Range("H27").Value = 1217.4286
$ws.Range("I27").Value = 982.625
$ws.Range("J27").Value = 1530.5
$ws.Range("K27").Value = 982.625
$ws.Range("L27").Value = 1530.5
$ws.Range("M27").Value = -875.625
$ws.Range("N27").Value = -1744.5
$ws.Range("H61").Value = 4261.619
$ws.Range("I61").Value = 4233.2
$ws.Range("K61").Value = 4233.2
$ws.Range("M61").Value = -4031.2
$ws.Range("H82").Value = 2065.75
$ws.Range("I82").Value = 1768.8182
$ws.Range("J82").Value = 2719
$ws.Range("K82").Value = 1768.8182
$ws.Range("L82").Value = 2719
$ws.Range("M82").Value = -1407.8182
$ws.Range("N82").Value = -3441
$ws.Range("H85").Value = 2065.75
$ws.Range("I85").Value = 1768.8182
$ws.Range("J85").Value = 2719
$ws.Range("K85").Value = 1768.8182
$ws.Range("L85").Value = 2719
$ws.Range("M85").Value = -520.8181999999999
$ws.Range("N85").Value = -5215
$ws.Range("H100").Value = 58886.156
$ws.Range("I100").Value = 63519.59
$ws.Range("K100").Value = 63519.59
$ws.Range("M100").Value = -62978.59
$ws.Range("H113").Value = 4261.619
$ws.Range("I113").Value = 4233.2
$ws.Range("K113").Value = 4233.2
$ws.Range("M113").Value = -2063.2
$ws.Range("H122").Value = 10572.728
$ws.Range("I122").Value = 8543
$ws.Range("J122").Value = 14124.75
$ws.Range("K122").Value = 25629
$ws.Range("L122").Value = 42374.25
$ws.Range("M122").Value = -23179
$ws.Range("N122").Value = -47274.25
$ws.Range("H126").Value = 15374.75
$ws.Range("I126").Value = 18499.666
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 55498.99800000001
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -53028.99800000001
$ws.Range("N126").Value = -22940
$ws.Range("H136").Value = 4615.737
$ws.Range("I136").Value = 4668.6875
$ws.Range("K136").Value = 14006.0625
$ws.Range("M136").Value = -11456.0625

# --- Sheet: WVR (15 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1008.5143
$ws.Range("I100").Value = 829.3077
$ws.Range("K100").Value = 1658.6154
$ws.Range("M100").Value = -1117.6154
$ws.Range("H132").Value = 3875.5483
$ws.Range("I132").Value = 3671.2334
$ws.Range("J132").Value = 10005
$ws.Range("K132").Value = 11013.7002
$ws.Range("L132").Value = 30015
$ws.Range("M132").Value = -8483.700199999999
$ws.Range("N132").Value = -35075
$ws.Range("H136").Value = 1868.3513
$ws.Range("I136").Value = 1660.5312
$ws.Range("K136").Value = 4981.5936
$ws.Range("M136").Value = -2431.5936
